$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 99.99995231628418
$ws.Range("D8").Value = 100
$ws.Range("D9").Value = 99.99996423721313
$ws.Range("D10").Value = 99.99994039535522
$ws.Range("D12").Value = 100
$ws.Range("D13").Value = 99.99994039535522
$ws.Range("D14").Value = 99.99998807907104
$ws.Range("D15").Value = 99.9112069606781
$ws.Range("D16").Value = 99.99693632125854
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 99.99977350234985
$ws.Range("D19").Value = 99.99994039535522
$ws.Range("D20").Value = 99.99997615814209
$ws.Range("D22").Value = 99.99998807907104
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 99.99994039535522
$ws.Range("D24").Value = 100
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = 99.9915599822998
$ws.Range("D33").Value = 100
$ws.Range("D34").Value = 100
$ws.Range("D36").Value = 100
$ws.Range("D38").Value = 99.99998807907104
$ws.Range("D39").Value = 99.99997615814209
$ws.Range("D41").Value = 99.99998807907104
$ws.Range("D42").Value = 99.99927282333374
$ws.Range("D44").Value = 99.99998807907104
$ws.Range("D45").Value = 99.99996423721313
$ws.Range("D47").Value = 99.9997615814209
$ws.Range("D48").Value = 99.90565180778503
$ws.Range("D49").Value = 99.9998927116394
$ws.Range("D50").Value = 99.99984502792358
$ws.Range("D51").Value = 99.99997615814209
$ws.Range("C52").Value = 2
$ws.Range("D52").Value = 99.99929666519165
$ws.Range("D53").Value = 99.99933242797852
$ws.Range("D55").Value = 100
$ws.Range("D56").Value = 99.99988079071045
$ws.Range("D57").Value = 99.99997615814209
$ws.Range("D59").Value = 99.99997615814209
$ws.Range("D62").Value = 100
$ws.Range("D63").Value = 89.74702954292297
$ws.Range("D65").Value = 100
$ws.Range("D70").Value = 99.99998807907104
$ws.Range("D71").Value = 99.99994039535522
$ws.Range("D72").Value = 100
$ws.Range("D74").Value = 100
$ws.Range("C78").Value = 1
$ws.Range("D78").Value = 88.20975422859192
$ws.Range("D81").Value = 99.9997615814209
$ws.Range("C82").Value = 2
$ws.Range("D82").Value = 99.99992847442627
$ws.Range("D83").Value = 100
$ws.Range("D85").Value = 99.99988079071045
$ws.Range("D87").Value = 99.99850988388062
$ws.Range("D88").Value = 82.39760994911194
$ws.Range("D89").Value = 74.03348088264465
$ws.Range("D90").Value = 99.99998807907104
$ws.Range("D92").Value = 99.98767375946045
$ws.Range("D93").Value = 100
$ws.Range("D94").Value = 99.42518472671509
$ws.Range("D95").Value = 99.99998807907104
$ws.Range("D98").Value = 99.89904761314392
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 99.91976618766785
$ws.Range("D100").Value = 99.99998807907104
$ws.Range("D101").Value = 99.99990463256836
$ws.Range("D102").Value = 99.99933242797852
$ws.Range("D104").Value = 99.99997615814209
$ws.Range("D105").Value = 99.99421834945679
$ws.Range("D106").Value = 100
$ws.Range("D108").Value = 99.98403787612915
$ws.Range("D109").Value = 100
$ws.Range("C112").Value = 2
$ws.Range("D112").Value = 99.99992847442627
$ws.Range("D113").Value = 99.99765157699585
$ws.Range("D115").Value = 99.99997615814209
$ws.Range("D117").Value = 100
$ws.Range("D118").Value = 100
$ws.Range("D120").Value = 100
$ws.Range("D121").Value = 99.99996423721313
$ws.Range("C123").Value = 1
$ws.Range("D123").Value = 100
$ws.Range("D124").Value = 99.99998807907104
$ws.Range("D126").Value = 100
$ws.Range("D127").Value = 86.68943643569946
$ws.Range("C128").Value = 1
$ws.Range("D128").Value = 95.32052278518677
$ws.Range("D129").Value = 99.99949932098389
$ws.Range("D130").Value = 99.99990463256836
$ws.Range("D131").Value = 100
$ws.Range("D137").Value = 99.9970555305481
$ws.Range("D138").Value = 99.99951124191284
$ws.Range("D139").Value = 100
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 50.00947117805481
$ws.Range("D141").Value = 99.99998807907104
$ws.Range("D142").Value = 99.99998807907104
$ws.Range("D146").Value = 100
$ws.Range("D148").Value = 99.99997615814209
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 59.84077453613281
$ws.Range("D150").Value = 99.99948740005493
$ws.Range("D151").Value = 99.94229078292847
$ws.Range("D152").Value = 99.99998807907104
$ws.Range("D153").Value = 99.99912977218628
$ws.Range("D154").Value = 100
$ws.Range("D155").Value = 100
$ws.Range("D156").Value = 100
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 99.99997615814209
